# Workbook: "Hortaliza, Vega Modelo de Temuco - Alcachofa"
# Two new daily price records were inserted into the data table, right
# before the existing row for date 44403 (2021-07-26). That pushes every
# subsequent record down by two rows, so the two oldest records that used
# to sit at the bottom of the sheet now land on two brand-new rows at the
# end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 132 - everything from 132 downward
# (previously ending at row 206) shifts down to end at row 208, matching
# the new <dimension ref="A1:R208"/>.
$ws.Range("A132:A133").EntireRow.Insert()

# New record -> row 132
$ws.Range("A132").Value = 10
$ws.Range("B132").Value = "Vega Modelo de Temuco"
$ws.Range("C132").Value = "La Araucanía"
$ws.Range("D132").Value = 44784
$ws.Range("E132").Value = 9
$ws.Range("F132").Value = 100112013
$ws.Range("G132").Value = "Alcachofa"
$ws.Range("H132").Value = "Española"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 350
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 20000
$ws.Range("M132").Value = 18857
$ws.Range("N132").Value = "`$/caja 30 unidades"
$ws.Range("O132").Value = "Provincia de Limarí"
$ws.Range("P132").Value = 629
$ws.Range("Q132").Value = 30
$ws.Range("R132").Value = "Hortaliza"

# New record -> row 133
$ws.Range("A133").Value = 10
$ws.Range("B133").Value = "Vega Modelo de Temuco"
$ws.Range("C133").Value = "La Araucanía"
$ws.Range("D133").Value = 44784
$ws.Range("E133").Value = 9
$ws.Range("F133").Value = 100112013
$ws.Range("G133").Value = "Alcachofa"
$ws.Range("H133").Value = "Madrigal"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 500
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = 15000
$ws.Range("N133").Value = "`$/caja 40 unidades"
$ws.Range("O133").Value = "Provincia de Limarí"
$ws.Range("P133").Value = 375
$ws.Range("Q133").Value = 40
$ws.Range("R133").Value = "Hortaliza"
